$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-151 was bumped by one day
# (serial date 45181 -> 45182, i.e. 2023-09-12 -> 2023-09-13).
for ($r = 2; $r -le 151; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
